$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 8; existing rows 8-41 shift down to 10-43.
$ws.Rows("8:9").Insert()

# New row 8: Early Diamond
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = [DateTime]"2022-01-04"
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100103
$ws.Cells.Item(8, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(8, 9).Value = 100103006
$ws.Cells.Item(8, 10).Value = "Nectarín"
$ws.Cells.Item(8, 11).Value = "Early Diamond"
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 250
$ws.Cells.Item(8, 14).Value = 19000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 19500
$ws.Cells.Item(8, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 1083
$ws.Cells.Item(8, 20).Value = 18

# New row 9: Super Queen
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value = [DateTime]"2022-01-04"
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100103
$ws.Cells.Item(9, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(9, 9).Value = 100103006
$ws.Cells.Item(9, 10).Value = "Nectarín"
$ws.Cells.Item(9, 11).Value = "Super Queen"
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 270
$ws.Cells.Item(9, 14).Value = 19000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 19500
$ws.Cells.Item(9, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 1083
$ws.Cells.Item(9, 20).Value = 18

Write-Host "done"
